# Update facultyCode values (column B) in the majors data sheet.
# The original file used a single shared string "GENERAL" (B2:B56) for the
# facultyCode column. This edit replaces it, for rows 3-56, with the actual
# faculty short-codes that correspond to each major, while row 2 (the
# "GENERAL" / general template row) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$facultyCodes = @{
    3  = "FBK"
    4  = "FIKF"
    5  = "FTI_DIPLOMA"
    6  = "FBK"
    7  = "FBK"
    8  = "FTI_DIPLOMA"
    9  = "PROFESI"
    10 = "PROFESI"
    11 = "PROFESI"
    12 = "FTI"
    13 = "FE"
    14 = "FE"
    15 = "FTSP"
    16 = "FTSP"
    17 = "FTSP"
    18 = "FE"
    19 = "FIKF"
    20 = "FIKOM"
    21 = "FIKOM"
    22 = "FTI"
    23 = "FTI"
    24 = "FIKF"
    25 = "FK"
    26 = "FE"
    27 = "FE"
    28 = "FSB"
    29 = "FPSI"
    30 = "FPSI"
    31 = "FSB"
    32 = "FSB"
    33 = "FIKTI"
    34 = "FIKTI"
    35 = "FIKTI"
    36 = "FTI"
    37 = "FTI"
    38 = "FTI"
    39 = "FTSP"
    40 = "FTSP"
    41 = "MAGISTER"
    42 = "MAGISTER"
    43 = "MAGISTER"
    44 = "MAGISTER"
    45 = "MAGISTER"
    46 = "MAGISTER"
    47 = "MAGISTER"
    48 = "MAGISTER"
    49 = "MAGISTER"
    50 = "MAGISTER"
    51 = "MAGISTER"
    52 = "PROFESI"
    53 = "PROFESI"
    54 = "PROFESI"
    55 = "PROFESI"
    56 = "PROFESI"
}

foreach ($row in $facultyCodes.Keys) {
    $ws.Cells.Item($row, 2).Value = $facultyCodes[$row]
}
